$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated price / volume(1h) values cell-by-cell, matching the
# authoritative diff. Numeric-looking strings (e.g. "217.44") are written
# with a leading apostrophe so Excel keeps them as text (matching the
# original inlineStr/text cell type) instead of silently converting them
# to numeric cells; the style is then reset to Normal so no stray
# quote-prefix / number-format is left behind on the cell.

$ws.Range("D2").Value = '27.174.92'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '1.648.76'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'217.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = "'0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.35%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = "'0.256"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.56%  '
$ws.Range("D9").Value = "'0.0628"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'19.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.10%  '
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").Value = '1.880.53'
$ws.Range("D13").Value = '1.694.29'
$ws.Range("E13").Value = '  +3.00%  '
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").Value = "'0.541"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.83%  '
$ws.Range("D16").Value = "'67.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.08%  '
$ws.Range("D17").Value = '27.186.82'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").Value = "'218.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D21").Value = "'6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.95%  '
$ws.Range("D22").Value = "'2.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.51%  '
$ws.Range("D23").Value = "'4.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("D24").Value = "'9.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("D25").Value = "'147.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("E26").Value = '  +2.95%  '
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = "'15.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("D30").Value = "'0.0509"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").Value = "'1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").Value = '  +1.15%  '
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("E34").Value = '  +1.76%  '
$ws.Range("D35").Value = '1.269.07'
$ws.Range("E35").Value = '  +1.95%  '
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("E38").Value = '  +3.20%  '
$ws.Range("D39").Value = "'0.849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").Value = "'5.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("E43").Value = '  +6.41%  '
$ws.Range("D44").Value = '1.791.01'
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("D45").Value = "'62.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.75%  '
$ws.Range("D46").Value = "'91.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("E48").Value = '  +1.88%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = "'0.0976"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = "'7.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.52%  '
